{"js": "  const replacements = [\n    [\"803\u00d79=7227\", \"367\u00d75=1835\"],\n    [\"953\u00d77=6671\", \"215\u00d72=430\"],\n    [\"270\u00d74=1080\", \"520\u00d79=4680\"],\n    [\"823\u00d73=2469\", \"816\u00d72=1632\"],\n    [\"349\u00d78=2792\", \"484\u00d73=1452\"],\n    [\"486\u00d76=2916\", \"850\u00d75=4250\"],\n    [\"932\u00d75=4660\", \"346\u00d72=692\"],\n    [\"104\u00d77=728\", \"665\u00d75=3325\"],\n    [\"584\u00d77=4088\", \"484\u00d79=4356\"],\n    [\"453\u00d75=2265\", \"969\u00d74=3876\"],\n    [\"270\u00d78=2160\", \"584\u00d73=1752\"],\n    [\"204\u00d78=1632\", \"107\u00d74=428\"],\n    [\"747\u00d75=3735\", \"362\u00d79=3258\"],\n    [\"166\u00d76=996\", \"154\u00d73=462\"],\n    [\"293\u00d78=2344\", \"130\u00d72=260\"],\n    [\"974\u00d75=4870\", \"788\u00d77=5516\"],\n    [\"734\u00d77=5138\", \"869\u00d79=7821\"],\n    [\"926\u00d79=8334\", \"629\u00d74=2516\"],\n    [\"526\u00d72=1052\", \"605\u00d77=4235\"],\n    [\"493\u00d79=4437\", \"834\u00d78=6672\"],\n    [\"531\u00d77=3717\", \"834\u00d73=2502\"],\n    [\"951\u00d77=6657\", \"142\u00d73=426\"],\n    [\"568\u00d79=5112\", \"849\u00d76=5094\"],\n    [\"934\u00d75=4670\", \"235\u00d77=1645\"],\n    [\"131\u00d75=655\", \"292\u00d79=2628\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      throw new Error(\"Could not find text: \" + oldText);\n    }\n\n    for (const range of results.items) {\n      range.insertText(newText, \"Replace\");\n    }\n  }\n\n  await context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"803\u00d79=7227\"; New = \"367\u00d75=1835\" },\n    @{ Old = \"953\u00d77=6671\"; New = \"215\u00d72=430\" },\n    @{ Old = \"270\u00d74=1080\"; New = \"520\u00d79=4680\" },\n    @{ Old = \"823\u00d73=2469\"; New = \"816\u00d72=1632\" },\n    @{ Old = \"349\u00d78=2792\"; New = \"484\u00d73=1452\" },\n    @{ Old = \"486\u00d76=2916\"; New = \"850\u00d75=4250\" },\n    @{ Old = \"932\u00d75=4660\"; New = \"346\u00d72=692\" },\n    @{ Old = \"104\u00d77=728\"; New = \"665\u00d75=3325\" },\n    @{ Old = \"584\u00d77=4088\"; New = \"484\u00d79=4356\" },\n    @{ Old = \"453\u00d75=2265\"; New = \"969\u00d74=3876\" },\n    @{ Old = \"270\u00d78=2160\"; New = \"584\u00d73=1752\" },\n    @{ Old = \"204\u00d78=1632\"; New = \"107\u00d74=428\" },\n    @{ Old = \"747\u00d75=3735\"; New = \"362\u00d79=3258\" },\n    @{ Old = \"166\u00d76=996\"; New = \"154\u00d73=462\" },\n    @{ Old = \"293\u00d78=2344\"; New = \"130\u00d72=260\" },\n    @{ Old = \"974\u00d75=4870\"; New = \"788\u00d77=5516\" },\n    @{ Old = \"734\u00d77=5138\"; New = \"869\u00d79=7821\" },\n    @{ Old = \"926\u00d79=8334\"; New = \"629\u00d74=2516\" },\n    @{ Old = \"526\u00d72=1052\"; New = \"605\u00d77=4235\" },\n    @{ Old = \"493\u00d79=4437\"; New = \"834\u00d78=6672\" },\n    @{ Old = \"531\u00d77=3717\"; New = \"834\u00d73=2502\" },\n    @{ Old = \"951\u00d77=6657\"; New = \"142\u00d73=426\" },\n    @{ Old = \"568\u00d79=5112\"; New = \"849\u00d76=5094\" },\n    @{ Old = \"934\u00d75=4670\"; New = \"235\u00d77=1645\" },\n    @{ Old = \"131\u00d75=655\"; New = \"292\u00d79=2628\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($pair.Old)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
